# nanocopter-partlist.xlsx - "all parts & values defined" edit
#
# Adds the still-missing part rows to the sheet:
#   - row 15: a reminder note about the Molex SMD connectors
#   - row 27: the MBR0520 sepic/charge-pump diode
#   - rows 30-44: the whole new "Capacitors" / "Inductors" / "Ferrite" /
#     "Widerstaende" block (Sepic circuit parts)
#
# The category/section headers (A30, A31, A37, A41, A43) reuse the same
# bold header style ("s=1") that is already used for every other
# category cell on the sheet (e.g. A18 "Spezialteile", A23 "Mosfets", ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Text values -------------------------------------------------------
# Written in the same order the parts were originally researched/typed in
# (scattered across the new rows rather than strictly top-to-bottom), so
# that cells which repeat the same text (e.g. BLM18KG331SN1D in B41/C41/E41)
# resolve to one shared value, just like the rest of the sheet already does.
$ws.Range("B32").Value = "C_IN 2.2uF Low ESR"
$ws.Range("E32").Value = "963-LMK212BJ225MG-T"
$ws.Range("C33").Value = "LMK212BJ105MG-T"
$ws.Range("B33").Value = "C_SW 1uF Low ESR"
$ws.Range("E33").Value = "963-LMK212BJ105MG-T"
$ws.Range("C35").Value = "JMK107BJ475MA-T"
$ws.Range("E35").Value = "963-JMK107BJ475MA-T"
$ws.Range("A37").Value = "Inductors"
$ws.Range("C37").Value = "IFSC1111AZER100M01"
$ws.Range("E37").Value = "70-IFSC1111AZER100M0"
$ws.Range("C38").Value = "LQH32CN100K33L"
$ws.Range("E38").Value = "81-LQH32CN100K33L"
$ws.Range("B38").Value = "Sepic: 10uH Murata"
$ws.Range("B37").Value = "Sepic: 10uH Vishay (alternative)"
$ws.Range("C39").Value = "CBC2518T100K"
$ws.Range("E39").Value = "963-CBC2518T100K"
$ws.Range("B39").Value = "Sepic: 10uH Taiyo Yuden (alternative 2) (KLEIN!)"
$ws.Range("B27").Value = "MBR0520 (SEPIC, ChargePump out)"
$ws.Range("B35").Value = "DCDC C_OUT & STM32 Bypass 4.7uF Low ESR"
$ws.Range("A41").Value = "Ferrite"
$ws.Range("B41").Value = "BLM18KG331SN1D"
$ws.Range("C41").Value = "BLM18KG331SN1D"
$ws.Range("E41").Value = "BLM18KG331SN1D"
$ws.Range("A43").Value = "Widerstände"
$ws.Range("B43").Value = "Spannungsteiler Sepic: 110k 1%"
$ws.Range("C43").Value = "CR0603-FX-1103GLF"
$ws.Range("B44").Value = "Spannungsteiler Sepic: 180k 1%"
$ws.Range("C44").Value = "CR0603-FX-1803ELF"
$ws.Range("B15").Value = "TODO: Richtige molex buchsen (smd) suchen"
$ws.Range("C27").Value = "MBR0520LT1G"
$ws.Range("E27").Value = "863-MBR0520LT1G"
$ws.Range("A30").Value = "Capacitors"
$ws.Range("A31").Value = "Sepic"
$ws.Range("C31").Value = "JMK316BJ106ML-T"
$ws.Range("E31").Value = "963-JMK316BJ106ML-T"
$ws.Range("B31").Value = "C_OUT 10uF Low ESR"
$ws.Range("C32").Value = "LMK212BJ225MG-T"

# --- 2) QTY / Preis numeric cells ----------------------------------------
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 0.35

$ws.Range("D31").Value = 1
$ws.Range("F31").Value = 0.46
$ws.Range("D32").Value = 1
$ws.Range("F32").Value = 0.24
$ws.Range("D33").Value = 1
$ws.Range("F33").Value = 0.36

$ws.Range("D35").Value = 2
$ws.Range("F35").Value = 0.27

$ws.Range("D37").Value = 2
$ws.Range("F37").Value = 0.53
$ws.Range("D38").Value = 2
$ws.Range("F38").Value = 0.28999999999999998
$ws.Range("D39").Value = 2
$ws.Range("F39").Value = 0.17

$ws.Range("D41").Value = 1
$ws.Range("F41").Value = 0.34

$ws.Range("D43").Value = 1
$ws.Range("D44").Value = 1

# --- 3) Bold the new category / section headers --------------------------
# (re-applies the workbook's existing bold header style, same as the other
# category cells like A18/A23/A26 already use)
$ws.Range("A30").Font.Bold = $true
$ws.Range("A31").Font.Bold = $true
$ws.Range("A37").Font.Bold = $true
$ws.Range("A41").Font.Bold = $true
$ws.Range("A43").Font.Bold = $true

# --- 4) View state ---------------------------------------------------------
# Scroll the window down a bit and leave the active selection on B16,
# matching where editing left off.
$excel.ActiveWindow.ScrollRow = 3
[void]$ws.Range("B16").Select()
